# Update crypto price/volume figures for Thu May  4 11:06:39 UTC 2023 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell: untouched, default ("General"/no explicit style) formatting,
# used to restore the original style index on cells where we must force
# text storage (NumberFormat "@") for values that would otherwise be
# auto-parsed as numbers (e.g. "1.004", "119.10") by Excel's input parser.
$blankStyle = $ws.Range("Z100").Style

$ws.Range("D2").Value = "29.268.34"
$ws.Range("E2").Value = "  +1.85%  "

$ws.Range("D3").Value = "1.914.64"
$ws.Range("E3").Value = "  +2.08%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = $blankStyle
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.82"
$ws.Range("D5").Style = $blankStyle
$ws.Range("E5").Value = "  +0.66%  "

$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4623"
$ws.Range("D7").Style = $blankStyle
$ws.Range("E7").Value = "  +0.59%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3948"
$ws.Range("D8").Style = $blankStyle
$ws.Range("E8").Value = "  +2.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.76"
$ws.Range("D9").Style = $blankStyle
$ws.Range("E9").Value = "  +1.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07948"
$ws.Range("D10").Style = $blankStyle
$ws.Range("E10").Value = "  +1.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.004"
$ws.Range("D11").Style = $blankStyle
$ws.Range("E11").Value = "  +0.69%  "

$ws.Range("E12").Value = "  +2.88%  "

$ws.Range("D13").Value = "1.934.77"
$ws.Range("E13").Value = "  +2.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.111"
$ws.Range("D14").Style = $blankStyle
$ws.Range("E14").Value = "  +1.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.771"
$ws.Range("D15").Style = $blankStyle
$ws.Range("E15").Value = "  +1.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06953"
$ws.Range("D16").Style = $blankStyle
$ws.Range("E16").Value = "  -0.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.61"
$ws.Range("D17").Style = $blankStyle
$ws.Range("E17").Value = "  +0.15%  "

$ws.Range("E18").Value = "  -0.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001011"
$ws.Range("D19").Style = $blankStyle
$ws.Range("E19").Value = "  +0.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.13"
$ws.Range("D20").Style = $blankStyle
$ws.Range("E20").Value = "  +1.57%  "

$ws.Range("E21").Value = "  -0.34%  "

$ws.Range("D22").Value = "29.298.55"
$ws.Range("E22").Value = "  +1.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.366"
$ws.Range("D23").Style = $blankStyle
$ws.Range("E23").Value = "  +1.65%  "

$ws.Range("E24").Value = "  +1.03%  "

$ws.Range("D25").Value = "2.140.52"
$ws.Range("E25").Value = "  +1.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.063"
$ws.Range("D26").Style = $blankStyle
$ws.Range("E26").Value = "  -2.97%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.79"
$ws.Range("D27").Style = $blankStyle
$ws.Range("E27").Value = "  +1.95%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.53"
$ws.Range("D28").Style = $blankStyle
$ws.Range("E28").Value = "  +1.53%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.102"
$ws.Range("D29").Style = $blankStyle
$ws.Range("E29").Value = "  +5.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.001"
$ws.Range("D30").Style = $blankStyle
$ws.Range("E30").Value = "  +1.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.10"
$ws.Range("D31").Style = $blankStyle

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09403"
$ws.Range("D32").Style = $blankStyle
$ws.Range("E32").Value = "  +0.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9290"
$ws.Range("D33").Style = $blankStyle
$ws.Range("E33").Value = "  +1.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.359"
$ws.Range("D34").Style = $blankStyle
$ws.Range("E34").Value = "  +1.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.360"
$ws.Range("D35").Style = $blankStyle
$ws.Range("E35").Value = "  +1.39%  "

$ws.Range("E36").Value = "  -1.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.207"
$ws.Range("D37").Style = $blankStyle
$ws.Range("E37").Value = "  +4.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05846"
$ws.Range("D38").Style = $blankStyle
$ws.Range("E38").Value = "  +1.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02109"
$ws.Range("D39").Style = $blankStyle
$ws.Range("E39").Value = "  +1.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.979"
$ws.Range("D40").Style = $blankStyle
$ws.Range("E40").Value = "  +3.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.002"
$ws.Range("D41").Style = $blankStyle
$ws.Range("E41").Value = "  -0.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5769"
$ws.Range("D42").Style = $blankStyle
$ws.Range("E42").Value = "  +2.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1807"
$ws.Range("D43").Style = $blankStyle
$ws.Range("E43").Value = "  +0.91%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.01"
$ws.Range("D44").Style = $blankStyle
$ws.Range("E44").Value = "  +0.97%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.284"
$ws.Range("D45").Style = $blankStyle
$ws.Range("E45").Value = "  +6.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.99"
$ws.Range("D46").Style = $blankStyle
$ws.Range("E46").Value = "  +1.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5431"
$ws.Range("D47").Style = $blankStyle
$ws.Range("E47").Value = "  +2.72%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.07077"
$ws.Range("D48").Style = $blankStyle
$ws.Range("E48").Value = "  -1.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.883"
$ws.Range("D49").Style = $blankStyle
$ws.Range("E49").Value = "  +3.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.561"
$ws.Range("D50").Style = $blankStyle
$ws.Range("E50").Value = "  +6.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.44"
$ws.Range("D51").Style = $blankStyle
$ws.Range("E51").Value = "  -0.19%  "
